# "added ERD and updated risk assess"
# Update the risk-assessment table text, row heights and column widths,
# add a new "Proposed Mitigation" entry for the cloud-hosting risk (J7),
# and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New proposed-mitigation cell for the "cloud may run into problems" risk row.
# Clone the formatting (centered + wrapped, same as the other Proposed
# Mitigation cells) from J4 before writing the new text into J7.
$ws.Range("J4").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J7").Value = "The app will be hosted on the cloud so therefore the cloud provider is responsible for the upkeep of the website"

# --- Update / extend description & mitigation text -------------------------
$ws.Range("C4").Value = "Unauthorized people may gain access to the website and view private information alter the data"
$ws.Range("C5").Value = "A hacker may try to gain access to the databse to steal information"
$ws.Range("J5").Value = "Keeps an eye for any attackers trying to gain access to the webapp and act accordingly"

# --- Row heights (rewrap after the longer / shorter text) ------------------
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 75

# --- Column widths -----------------------------------------------------
# Column C: widen to fit the longer description text (drop autofit/bestFit).
$ws.Columns.Item(3).ColumnWidth = 25.8
# Column J: new explicit width for the Proposed Mitigation column.
$ws.Columns.Item(10).ColumnWidth = 16.1

# --- Selection ---------------------------------------------------------
$ws.Range("C8").Select()
